# "mais um update, a partir de agora e sempre a abrir"
#
# The task list got a new "% de Completude" column, a handful of new
# tasks (Toasts / MyApplication helpers / Spinner Cozinha), and the old
# "Data de realização" column was dropped. Rebuild B3:D28 from scratch
# with the final layout, then restyle/resize to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# --- Row data -------------------------------------------------------
# Each entry: row number, B (task), C (priority), D (% complete).
# $null means "leave blank".
$rows = @(
    @{ R = 3;  B = "Lista de Tarefas a Realizar"; C = "Grau de Prioridade (1 a 5, 5 máximo)"; D = "% de Completude" }
    @{ R = 4;  B = "Meter prints dos acontecimentos no servidor"; C = 2; D = 10 }
    @{ R = 5;  B = "Corrigir layouts"; C = 3; D = 100 }
    @{ R = 6;  B = "Quando aplicação abre, ir buscar estados ao servidor"; C = 5; D = 40 }
    @{ R = 7;  B = "Quarto, guardar modos no servidor e respetivas mudanças"; C = 4; D = 0 }
    @{ R = 8;  B = "Passar ar condicionados para o servidor"; C = 3; D = 0 }
    @{ R = 9;  B = "Passar estado da televisão para o servidor"; C = 3; D = 0 }
    @{ R = 10; B = "Passar águas da casa de banho para o servidor"; C = 3; D = 100 }
    @{ R = 11; B = "Passar estados do forno para o servidor"; C = 3; D = 0 }
    @{ R = 12; B = "Passar estados do microondas para o servidor"; C = 3; D = 0 }
    @{ R = 13; B = "Arranjar quantro pontos na sala 121 onde a aplicação muda para quatro localizações diferentes"; C = 5; D = 30 }
    @{ R = 14; B = "Completar comando de televisão"; C = 1; D = 0 }
    @{ R = 15; B = "Guardar a lista de canais mais visots"; C = 1; D = 0 }
    @{ R = 16; B = "Feedback de operações na aplicação movel"; C = 2; D = 30 }
    @{ R = 17; B = "Problemas de Toasts"; C = 4; D = 0 }
    @{ R = 18; B = "Mudar Class do MyApplication para usar os helpers"; C = 5; D = 80 }
    @{ R = 19; B = "Completar Helpers"; C = 5; D = 60 }
    @{ R = 20; B = "Spinner Cozinha"; C = 1; D = 0 }
    @{ R = 21; B = $null; C = $null; D = $null }
    @{ R = 22; B = "Testes"; C = $null; D = $null }
    @{ R = 23; B = "Testar se o servidor fica em modo manual ou automático"; C = 5; D = 0 }
    @{ R = 24; B = "Testar mudanças de estado da luz e janelas no modo manual"; C = 5; D = 0 }
    @{ R = 25; B = "Testar mudanças de dia/noite e respectivas mudanças nos comportamentos de janelas e luz"; C = 5; D = 0 }
    @{ R = 26; B = "Outros testes (colocar aqui)"; C = $null; D = $null }
    @{ R = 27; B = $null; C = $null; D = $null }
    @{ R = 28; B = "Reportar erros de testes e faltas de coisas a fazer na aplicação"; C = $null; D = $null }
)

# --- Wipe the old B:D block first (old sheet only went to row 24, a
# couple of rows change from "has C" to "no C", so clear before rewrite) --
$ws.Range("B3:D28").ClearContents() | Out-Null

# --- Write values -----------------------------------------------------
foreach ($row in $rows) {
    $r = $row.R
    if ($null -ne $row.B) { $ws.Cells.Item($r, 2).Value = $row.B }
    if ($null -ne $row.C) { $ws.Cells.Item($r, 3).Value = $row.C }
    if ($null -ne $row.D) { $ws.Cells.Item($r, 4).Value = $row.D }
}

# --- Header row (bold) --------------------------------------------------
$ws.Range("B3:D3").Font.Bold = $true

# --- Section headers / footer (bold, same style as the title row) ------
$ws.Range("B22").Font.Bold = $true
$ws.Range("B28").Font.Bold = $true

# --- Center-align the Priority / % Completude data columns -------------
$ws.Range("C4:C21").HorizontalAlignment = $xlCenter
$ws.Range("D4:D21").HorizontalAlignment = $xlCenter
$ws.Range("C23:C25").HorizontalAlignment = $xlCenter
$ws.Range("D23:D25").HorizontalAlignment = $xlCenter
$ws.Range("D3").HorizontalAlignment = $xlCenter

# --- Column widths (bestFit-style, matching the final C/D columns) -----
$ws.Columns("C").ColumnWidth = 31.3
$ws.Columns("D").ColumnWidth = 14.5

# --- Selection left where the author's last edit landed ----------------
$ws.Range("F10:F11").Select() | Out-Null
